{"js": "// Translate the English labels/values in the facilitation guide table to\n// Amharic. Each entry below is a unique (englishText -> amharicText) pair\n// found in the document; we use Range.search with matchCase to locate the\n// exact run and then replace its text in place (preserving formatting).\nconst replacements = [\n  [\"Video Title\", \"\u12e8\u126a\u12f2\u12ee \u12a0\u122d\u12d5\u1235\u1275\"],\n  [\"A Lesson on Conway\\u2019s Game of Life [Samuel Okoth]\", \"\u12e8\u12ae\u1295\u12ce\u12ed \u12e8\u1205\u12ed\u12c8\u1275 \u132b\u12c8\u1273[\u1233\u1219\u12a4\u120d \u12a6\u12ae\u12dd]\"],\n  [\"Topic\", \"\u12e8\u1275\u121d\u1205\u122d\u1275 \u122d\u12d5\u1235\"],\n  [\"Mathematical Thinking\", \"\u1212\u1233\u1263\u12ca \u12a5\u1233\u1264\"],\n  [\"Aim(s)\", \"\u12a0\u120b\u121b(\u12ce\u127d)\"],\n  [\n    \"Define an interesting and unpredictable cell automaton. For example, discover some configurations that last for a long time before dying and other configurations to go on forever without allowing cycles.\",\n    \"\u1233\u1262 \u12a5\u1293 \u12e8\u121b\u12ed\u1308\u1218\u1275 \u12e8\u1205\u12cb\u1235 \u12a0\u12cd\u1276\u121b\u1276\u1295 \u1218\u134d\u1320\u122d\u1361\u1361 \u1218\u134d\u1320\u122d\u1361\u1361 \u1208\u121d\u1233\u120c\u1363 \u12e8\u1270\u12c8\u1230\u1291 \u12cd\u1245\u122e\u127d\u1295 \u1218\u1348\u1208\u130d \u1208\u1228\u1305\u121d \u130a\u12dc \u1233\u12ed\u121e\u1271 \u12e8\u121a\u1246\u12e9 \u12a0\u1293 \u120c\u120b \u12cd\u1245\u122d \u1208\u12d8\u120b\u1208\u121d \u12d1\u12f0\u1275\u1295 \u1233\u12ed\u1348\u1245\u12f5 \u12e8\u121a\u1204\u12f5\",\n  ],\n  [\"Length\", \"\u122d\u12dd\u1218\u1275\"],\n  [\"Camp Location\", \"\u12e8\u12ab\u121d\u1355 \u12a0\u12f5\u122b\u123b\"],\n  [\"Facilitators\", \"\u12a0\u1235\u1270\u1263\u1263\u122a\u12ce\u127d\"],\n  [\"N. of students\", \"\u12e8 \u1270\u121b\u122a\u12ce\u127d \u1265\u12db\u1275\"],\n  [\"Date\", \"\u1240\u1295\"],\n  [\"Resources\", \"\u12a0\u1235\u1348\u120b\u130a\"],\n  [\"needed\", \"\u1218\u1233\u122a\u12eb\u12ce\u127d\"],\n  [\"Paper to draw square grids, 2 different coloured post-its\", \"\u12c8\u1228\u1240\u1275 \u12ab\u122c \u1348\u122d\u130d\u122d\u130d \u1208\u1218\u1235\u122b\u1275\u1363 2 \u12e8\u1270\u1208\u12eb\u12e9 \u1263\u1208\u1240\u1208\u121d \u120d\u1325\u134e\u127d\"],\n  [\"Preparations\", \"\u12e0\u130d\u1301\u1290\u1276\u127d\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + from);\n  }\n\n  results.items[0].insertText(to, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Translate the English labels/values in the facilitation guide table to\n# Amharic. Each pair below is a unique (englishText -> amharicText) mapping\n# found in the document; we use Range.Find/Execute (case-sensitive, so that\n# e.g. \"Facilitators\" doesn't also match the lower-case \"facilitators\" that\n# appears inside unrelated placeholder text) to locate the exact text and\n# replace it in place, preserving the run's formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Video Title\", \"\u12e8\u126a\u12f2\u12ee \u12a0\u122d\u12d5\u1235\u1275\"),\n    @(\"A Lesson on Conway\u2019s Game of Life [Samuel Okoth]\", \"\u12e8\u12ae\u1295\u12ce\u12ed \u12e8\u1205\u12ed\u12c8\u1275 \u132b\u12c8\u1273[\u1233\u1219\u12a4\u120d \u12a6\u12ae\u12dd]\"),\n    @(\"Topic\", \"\u12e8\u1275\u121d\u1205\u122d\u1275 \u122d\u12d5\u1235\"),\n    @(\"Mathematical Thinking\", \"\u1212\u1233\u1263\u12ca \u12a5\u1233\u1264\"),\n    @(\"Aim(s)\", \"\u12a0\u120b\u121b(\u12ce\u127d)\"),\n    @(\"Define an interesting and unpredictable cell automaton. For example, discover some configurations that last for a long time before dying and other configurations to go on forever without allowing cycles.\",\n      \"\u1233\u1262 \u12a5\u1293 \u12e8\u121b\u12ed\u1308\u1218\u1275 \u12e8\u1205\u12cb\u1235 \u12a0\u12cd\u1276\u121b\u1276\u1295 \u1218\u134d\u1320\u122d\u1361\u1361 \u1218\u134d\u1320\u122d\u1361\u1361 \u1208\u121d\u1233\u120c\u1363 \u12e8\u1270\u12c8\u1230\u1291 \u12cd\u1245\u122e\u127d\u1295 \u1218\u1348\u1208\u130d \u1208\u1228\u1305\u121d \u130a\u12dc \u1233\u12ed\u121e\u1271 \u12e8\u121a\u1246\u12e9 \u12a0\u1293 \u120c\u120b \u12cd\u1245\u122d \u1208\u12d8\u120b\u1208\u121d \u12d1\u12f0\u1275\u1295 \u1233\u12ed\u1348\u1245\u12f5 \u12e8\u121a\u1204\u12f5\"),\n    @(\"Length\", \"\u122d\u12dd\u1218\u1275\"),\n    @(\"Camp Location\", \"\u12e8\u12ab\u121d\u1355 \u12a0\u12f5\u122b\u123b\"),\n    @(\"Facilitators\", \"\u12a0\u1235\u1270\u1263\u1263\u122a\u12ce\u127d\"),\n    @(\"N. of students\", \"\u12e8 \u1270\u121b\u122a\u12ce\u127d \u1265\u12db\u1275\"),\n    @(\"Date\", \"\u1240\u1295\"),\n    @(\"Resources\", \"\u12a0\u1235\u1348\u120b\u130a\"),\n    @(\"needed\", \"\u1218\u1233\u122a\u12eb\u12ce\u127d\"),\n    @(\"Paper to draw square grids, 2 different coloured post-its\", \"\u12c8\u1228\u1240\u1275 \u12ab\u122c \u1348\u122d\u130d\u122d\u130d \u1208\u1218\u1235\u122b\u1275\u1363 2 \u12e8\u1270\u1208\u12eb\u12e9 \u1263\u1208\u1240\u1208\u121d \u120d\u1325\u134e\u127d\"),\n    @(\"Preparations\", \"\u12e0\u130d\u1301\u1290\u1276\u127d\")\n)\n\nforeach ($pair in $replacements) {\n    $source = $pair[0]\n    $target = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($source, $true, $false, $false, $false, $false, $true, 1, $false, $target, 2)\n}\n"}
